$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "72.163.90"
$ws.Range("E2").Value = "  +0.24%  "

Set-TextValue "D3" "4.040.90"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("E4").Value = "  -0.10%  "

Set-TextValue "D5" "538.74"
$ws.Range("E5").Value = "  +0.33%  "

Set-TextValue "D6" "150.84"
$ws.Range("E6").Value = "  -1.82%  "

Set-TextValue "D7" "4.034.03"
$ws.Range("E7").Value = "  +0.00%  "

Set-TextValue "D8" "0.701"
$ws.Range("E8").Value = "  +0.69%  "

Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.06%  "

Set-TextValue "D10" "0.753"
$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("E11").Value = "  -1.40%  "

Set-TextValue "D12" "53.81"
$ws.Range("E12").Value = "  +10.37%  "

$ws.Range("E13").Value = "  -1.21%  "

Set-TextValue "D14" "10.88"
$ws.Range("E14").Value = "  -0.67%  "

Set-TextValue "D15" "4.677.66"
$ws.Range("E15").Value = "  -0.15%  "

Set-TextValue "D16" "4.034.61"
$ws.Range("E16").Value = "  -0.40%  "

Set-TextValue "D17" "14.37"
$ws.Range("E17").Value = "  -0.25%  "

Set-TextValue "D18" "20.58"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("E20").Value = "  -0.91%  "

Set-TextValue "D21" "72.104.24"
$ws.Range("E21").Value = "  +0.19%  "

Set-TextValue "D22" "442.57"
$ws.Range("E22").Value = "  +1.41%  "

Set-TextValue "D23" "97.14"
$ws.Range("E23").Value = "  -2.46%  "

Set-TextValue "D24" "3.51"
$ws.Range("E24").Value = "  -2.13%  "

Set-TextValue "D25" "14.67"
$ws.Range("E25").Value = "  -0.30%  "

Set-TextValue "D26" "4.24"
$ws.Range("E26").Value = "  -0.34%  "

Set-TextValue "D27" "4.31"
$ws.Range("E27").Value = "  +16.18%  "

Set-TextValue "D28" "11.25"
$ws.Range("E28").Value = "  -0.40%  "

Set-TextValue "D29" "10.77"
$ws.Range("E29").Value = "  -2.16%  "

Set-TextValue "D30" "5.95"
$ws.Range("E30").Value = "  +2.00%  "

Set-TextValue "D31" "37.14"
$ws.Range("E31").Value = "  -0.11%  "

Set-TextValue "D32" "8.21"
$ws.Range("E32").Value = "  +17.96%  "

$ws.Range("E33").Value = "  +1.12%  "

Set-TextValue "D34" "13.54"
$ws.Range("E34").Value = "  -1.57%  "

Set-TextValue "D35" "49.32"
$ws.Range("E35").Value = "  +14.55%  "

Set-TextValue "D36" "682.46"
$ws.Range("E36").Value = "  +0.01%  "

Set-TextValue "D37" "66.90"
$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("E38").Value = "  +4.21%  "

Set-TextValue "D39" "0.0₃0871"
$ws.Range("E39").Value = "  +2.65%  "

Set-TextValue "D40" "0.149"
$ws.Range("E40").Value = "  -5.83%  "

Set-TextValue "D41" "3.39"
$ws.Range("E41").Value = "  -2.14%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D42" "11.21"
$ws.Range("E42").Value = "  +17.47%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D43" "3.36"
$ws.Range("E43").Value = "  -4.18%  "

$ws.Range("E44").Value = "  +0.21%  "

Set-TextValue "D45" "0.0494"
$ws.Range("E45").Value = "  -0.35%  "

Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.08%  "

Set-TextValue "D47" "0.151"
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("E48").Value = "  -3.43%  "

Set-TextValue "D49" "3.10"
$ws.Range("E49").Value = "  +0.81%  "

Set-TextValue "D50" "3.30"
$ws.Range("E50").Value = "  -3.68%  "

Set-TextValue "D51" "3.42"
$ws.Range("E51").Value = "  +1.44%  "
